{"js": "// The template was re-generated, which re-minted the GUID used to tag the\n// \"REF testid\" field's runs (w:rsidR) and the internal id of the\n// \"testid\" bookmark that follows it. Reproduce that by rewriting the\n// paragraph that holds the field + bookmarkStart with the new GUID baked\n// into every run, via a literal OOXML splice (so the exact attribute\n// value is preserved, not just \"some new run\").\n\nconst NEW_RUN_GUID = \"016BA9042AD3B3B853890BD84BD880C4\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find the paragraph that contains the \" REF testid \\h \" field (it is the\n// paragraph whose text resolves to the field result \"testidref\").\nlet target = null;\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"testidref\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  const ooxml =\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    '<w:p w:rsidR=\"00FA23FF\" w:rsidRPr=\"00DC5685\" w:rsidRDefault=\"00FA23FF\" w:rsidP=\"00FA23FF\">' +\n    \"<w:pPr><w:rPr><w:lang w:val=\\\"en-US\\\"/></w:rPr></w:pPr>\" +\n    \"<w:r/>\" +\n    '<w:r w:rsidR=\"' + NEW_RUN_GUID + '\"><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RUN_GUID + '\"><w:instrText xml:space=\"preserve\"> REF testid \\\\h </w:instrText></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RUN_GUID + '\"><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RUN_GUID + '\"><w:rPr><w:b w:val=\"true\"/><w:noProof/></w:rPr><w:t>testidref</w:t></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RUN_GUID + '\"><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n    '<w:bookmarkStart w:name=\"testid\" w:id=\"151181435939870919812341496777861426640\"/>' +\n    \"<w:r/>\" +\n    \"</w:p>\" +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\";\n\n  target.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The template was re-generated, which re-minted the GUID used to tag the\n# \"REF testid\" field's runs (w:rsidR) and the internal id of the\n# \"testid\" bookmark that follows it. Reproduce that by rewriting the\n# paragraph that holds the field + bookmarkStart with the new GUID baked\n# into every run, via a literal WordOpenXML splice (so the exact\n# attribute value is preserved, not just \"some new run\").\n\n$word = New-Object -ComObject Word.Application\n$d = $word.ActiveDocument\n\n$NEW_RUN_GUID = \"016BA9042AD3B3B853890BD84BD880C4\"\n\n# Find the paragraph that contains the \" REF testid \\h \" field (it is the\n# paragraph whose text resolves to the field result \"testidref\").\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text -like \"*testidref*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" w:rsidR=\"00FA23FF\" w:rsidRPr=\"00DC5685\" w:rsidRDefault=\"00FA23FF\" w:rsidP=\"00FA23FF\">' +\n        '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n        '<w:r/>' +\n        '<w:r w:rsidR=\"' + $NEW_RUN_GUID + '\"><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n        '<w:r w:rsidR=\"' + $NEW_RUN_GUID + '\"><w:instrText xml:space=\"preserve\"> REF testid \\h </w:instrText></w:r>' +\n        '<w:r w:rsidR=\"' + $NEW_RUN_GUID + '\"><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n        '<w:r w:rsidR=\"' + $NEW_RUN_GUID + '\"><w:rPr><w:b w:val=\"true\"/><w:noProof/></w:rPr><w:t>testidref</w:t></w:r>' +\n        '<w:r w:rsidR=\"' + $NEW_RUN_GUID + '\"><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n        '<w:bookmarkStart w:name=\"testid\" w:id=\"151181435939870919812341496777861426640\"/>' +\n        '<w:r/>' +\n        '</w:p>'\n\n    $target.Range.InsertXML($xml)\n}\n"}
